$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 1") {
        $p.Style = "Heading 2"
        $p.SpaceBefore = 2.0
        $p.SpaceAfter = 0
        $p.SpaceBeforeAuto = 0
        $p.SpaceAfterAuto = 0
        $p.LineSpacingRule = 5
        $p.LineSpacing = 12.95
        $p.Alignment = 0
        $p.LeftIndent = 0
        $p.RightIndent = 0
    }
}
